$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '68.688.05'
$ws.Range("E2").Value = '  +1.33%  '
$ws.Range("D3").Value = '3.819.12'
$ws.Range("E3").Value = '  +0.21%  '
$ws.Range("E4").Value = '  +0.33%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '613.31'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.66%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '164.82'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.74%  '
$ws.Range("D7").Value = '3.817.91'
$ws.Range("E7").Value = '  +0.29%  '
$ws.Range("E8").Value = '  +0.00%  '
$ws.Range("E9").Value = '  -0.08%  '
$ws.Range("E10").Value = '  +0.76%  '
$ws.Range("E11").Value = '  -0.20%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '6.72'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +5.68%  '
$ws.Range("E13").Value = '  -0.59%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '35.46'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.36%  '
$ws.Range("D15").Value = '4.461.72'
$ws.Range("E15").Value = '  +0.19%  '
$ws.Range("D16").Value = '3.826.51'
$ws.Range("E16").Value = '  +0.81%  '
$ws.Range("D17").Value = '68.665.67'
$ws.Range("E17").Value = '  +1.25%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '18.13'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.25%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.13'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.76%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '464.95'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.17%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.65'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.76%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.701'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.07%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.0000149'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +2.46%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '83.81'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.56%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '12.03'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.76%  '
$ws.Range("E27").Value = '  +0.16%  '
$ws.Range("E28").Value = '  +0.06%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '10.00'
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").Value = '3.966.07'
$ws.Range("E30").Value = '  +0.12%  '
$ws.Range("B31").Value = 'PancakeSwap'
$ws.Range("C31").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.64'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -5.31%  '
$ws.Range("B32").Value = 'ImmutableX'
$ws.Range("C32").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.22'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.13%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '7.26'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.95%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '29.04'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.35%  '
$ws.Range("E35").Value = '  -0.06%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '9.07'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.04%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.101'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.93%  '
$ws.Range("E38").Value = '  +6.43%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.90'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.74%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.979'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.87%  '
$ws.Range("E41").Value = '  -1.96%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.00'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.04%  '
$ws.Range("E43").Value = '  +0.04%  '
$ws.Range("B44").Value = 'TheGraph'
$ws.Range("C44").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.298'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.27%  '
$ws.Range("B45").Value = 'Monero'
$ws.Range("C45").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '153.68'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.27%  '
$ws.Range("B46").Value = 'OKB'
$ws.Range("C46").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '46.57'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -2.39%  '
$ws.Range("B47").Value = 'ONDO'
$ws.Range("C47").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.40'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.26%  '
$ws.Range("B48").Value = 'Arweave'
$ws.Range("C48").Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '42.65'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -4.42%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.38'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.37%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.88'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.51%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '378.42'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -3.07%  '
